# ManageProducts.xlsx - "Addition of test xml for jenkins"
#
# Three existing product rows (B2, B5, B8 on the "Input" sheet) are
# refreshed with newly generated product SKUs. Each updated cell also
# picks up the same "new product" look the sheet already uses elsewhere
# on the product column: a thin top/bottom border plus the sheet's
# standard shaded (indexed color 9) fill, so the freshly written rows are
# visually consistent with the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ProductSku {
    param(
        [string]$Address,
        [string]$Sku
    )

    $cell = $ws.Range($Address)
    $cell.Value = $Sku

    # Match the look-and-feel already used for product rows in this sheet.
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(8).Weight = 2      # xlThin
    $cell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $cell.Borders.Item(9).Weight = 2      # xlThin
    $cell.Interior.Pattern = 1            # xlSolid
    $cell.Interior.ColorIndex = 9
}

Set-ProductSku "B2" "produOai"
Set-ProductSku "B5" "prodzqOe"
Set-ProductSku "B8" "prodrHEk"
